# Updates Leve profit figures across several sheets (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H2").Value = 100
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -9540

$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -10872

$ws.Range("H80").Value = 3919.7
$ws.Range("I80").Value = 4099.8
$ws.Range("K80").Value = 12299.4
$ws.Range("M80").Value = -11301.4

$ws.Range("H83").Value = 3919.7
$ws.Range("I83").Value = 4099.8
$ws.Range("K83").Value = 36898.2
$ws.Range("M83").Value = -31906.2

$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 6050
$ws.Range("J88").Value = 4300
$ws.Range("K88").Value = 6050
$ws.Range("L88").Value = 4300
$ws.Range("M88").Value = -5644
$ws.Range("N88").Value = -5112

$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 6050
$ws.Range("J91").Value = 4300
$ws.Range("K91").Value = 6050
$ws.Range("L91").Value = 4300
$ws.Range("M91").Value = -4646
$ws.Range("N91").Value = -7108

$ws.Range("H113").Value = 8166.5
$ws.Range("I113").Value = 8199.799999999999
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 8199.799999999999
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -4945.799999999999
$ws.Range("N113").Value = -14508

$ws.Range("H131").Value = 3539
$ws.Range("J131").Value = 3225
$ws.Range("L131").Value = 9675
$ws.Range("N131").Value = -19755

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H45").Value = 2459.4
$ws.Range("I45").Value = 3063.8572
$ws.Range("J45").Value = 1049
$ws.Range("K45").Value = 3063.8572
$ws.Range("L45").Value = 1049
$ws.Range("M45").Value = -2686.8572
$ws.Range("N45").Value = -1803

$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H86").Value = 25000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 25000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 25000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -27246

$ws.Range("H89").Value = 25000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 25000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 125000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -136232

$ws.Range("H94").Value = 4510.75
$ws.Range("I94").Value = 4113
$ws.Range("J94").Value = 6499.5
$ws.Range("K94").Value = 4113
$ws.Range("L94").Value = 6499.5
$ws.Range("M94").Value = -3662
$ws.Range("N94").Value = -7401.5

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 1082.4166
$ws.Range("I16").Value = 1251.7142
$ws.Range("J16").Value = 845.4
$ws.Range("K16").Value = 1251.7142
$ws.Range("L16").Value = 845.4
$ws.Range("M16").Value = -964.7141999999999
$ws.Range("N16").Value = -1419.4

$ws.Range("H31").Value = 2149.5
$ws.Range("J31").Value = 2832.5
$ws.Range("L31").Value = 2832.5
$ws.Range("N31").Value = -3422.5

$ws.Range("H34").Value = 2149.5
$ws.Range("J34").Value = 2832.5
$ws.Range("L34").Value = 2832.5
$ws.Range("N34").Value = -3236.5

$ws.Range("H93").Value = 6400
$ws.Range("I93").Value = 6400
$ws.Range("K93").Value = 6400
$ws.Range("M93").Value = -4528

$ws.Range("H99").Value = 5526.875
$ws.Range("I99").Value = 5526.875
$ws.Range("K99").Value = 5526.875
$ws.Range("M99").Value = -4028.875

$ws.Range("H107").Value = 1138
$ws.Range("I107").Value = 1363.5714
$ws.Range("K107").Value = 1363.5714
$ws.Range("M107").Value = 556.4286

$ws.Range("H113").Value = 1082.4166
$ws.Range("I113").Value = 1251.7142
$ws.Range("J113").Value = 845.4
$ws.Range("K113").Value = 1251.7142
$ws.Range("L113").Value = 845.4
$ws.Range("M113").Value = 918.2858000000001
$ws.Range("N113").Value = -5185.4

$ws.Range("H122").Value = 4009.6667
$ws.Range("I122").Value = 4009.6667
$ws.Range("K122").Value = 12029.0001
$ws.Range("M122").Value = -9579.000100000001

$ws.Range("H126").Value = 5526.875
$ws.Range("I126").Value = 5526.875
$ws.Range("K126").Value = 16580.625
$ws.Range("M126").Value = -14110.625

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 1199.7142
$ws.Range("I5").Value = 1951.5
$ws.Range("J5").Value = 899
$ws.Range("K5").Value = 5854.5
$ws.Range("L5").Value = 2697
$ws.Range("M5").Value = -5742.5
$ws.Range("N5").Value = -2921

$ws.Range("H92").Value = 417.6
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 362.66666
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 1087.99998
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -3583.99998

$ws.Range("H97").Value = 559.7917
$ws.Range("I97").Value = 528.55
$ws.Range("J97").Value = 716
$ws.Range("K97").Value = 1585.65
$ws.Range("L97").Value = 2148
$ws.Range("M97").Value = -1089.65
$ws.Range("N97").Value = -3140

$ws.Range("H122").Value = 3460.7144
$ws.Range("I122").Value = 1248.5
$ws.Range("J122").Value = 3693.5789
$ws.Range("K122").Value = 11236.5
$ws.Range("L122").Value = 33242.2101
$ws.Range("M122").Value = -8786.5
$ws.Range("N122").Value = -38142.2101

$ws.Range("H135").Value = 1199.7142
$ws.Range("I135").Value = 1951.5
$ws.Range("J135").Value = 899
$ws.Range("K135").Value = 17563.5
$ws.Range("L135").Value = 8091
$ws.Range("M135").Value = -15028.5
$ws.Range("N135").Value = -13161

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H55").Value = 3000
$ws.Range("I55").Value = 3000
$ws.Range("K55").Value = 3000
$ws.Range("M55").Value = -2673

$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 3600.4
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 10801.2
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -8271.200000000001
$ws.Range("N132").Value = -20057

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H82").Value = 3059.375
$ws.Range("I82").Value = 3155.7144
$ws.Range("K82").Value = 3155.7144
$ws.Range("M82").Value = -2794.7144

$ws.Range("H85").Value = 3059.375
$ws.Range("I85").Value = 3155.7144
$ws.Range("K85").Value = 3155.7144
$ws.Range("M85").Value = -1907.7144

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H122").Value = 5888.5
$ws.Range("I122").Value = 5888.5
$ws.Range("K122").Value = 17665.5
$ws.Range("M122").Value = -15215.5
